$d = $word.ActiveDocument

# --- Locate the two empty trailing paragraphs that receive the new
# "Common pitfalls" / answer text. They are the 2nd and 4th paragraphs
# after the last paragraph that already has text ("The heavy math loops
# seem to work best ..."), i.e.:
#   anchor (text) -> empty -> "Common pitfalls:" -> empty -> "None known..."
$count = $d.Paragraphs.Count
$anchor = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*heavy math loops seem to work best*") {
        $anchor = $i
    }
}
$targetCommon = $anchor + 2
$targetAnswer = $anchor + 4

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'
$rPrBlock = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$runRPrBlock = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr>'

# --- Paragraph 1: "Common pitfalls:" -------------------------------------
# Keeps its existing <w:ind w:left="0" w:firstLine="0"/>; the run gains
# rFonts/sz/szCs and the new text.
$pCommon = $d.Paragraphs.Item($targetCommon)
$commonText = "Common pitfalls:"
$xmlCommon = "<w:p $wNs $w14Ns " + `
    'w:rsidR="00000000" w:rsidDel="00000000" w:rsidP="00000000" ' + `
    'w:rsidRDefault="00000000" w:rsidRPr="00000000" w14:paraId="00000064">' + `
    "<w:pPr>$rPrBlock</w:pPr>" + `
    '<w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000">' + `
    "$runRPrBlock" + `
    '<w:t xml:space="preserve">' + $commonText + '</w:t>' + `
    '</w:r>' + `
    '</w:p>'
$pCommon.Range.InsertXML($xmlCommon)

# Re-establish the paragraph's zero indent: InsertXML drops explicit
# zero-valued w:ind attributes, so set it again via the object model to
# materialize <w:ind w:left="0" w:firstLine="0"/>.
$pCommon2 = $d.Paragraphs.Item($targetCommon)
$pCommon2.Range.ParagraphFormat.LeftIndent = 0
$pCommon2.Range.ParagraphFormat.FirstLineIndent = 0

# --- Paragraph 2: the "None known ..." answer -----------------------------
# This one must additionally lose its <w:ind left="0" firstLine="0"/> for
# good (InsertXML's zero-stripping is exactly what we want here, so no
# indent is re-applied afterwards).
$pAnswer = $d.Paragraphs.Item($targetAnswer)
$answerText = "None known at this time, although it is expected that computer " + `
    "science students will have a difficult time with the theory and " + `
    "function of the N-body code.  All students might need some help with " + `
    "the Monte Carlo section."

$xmlAnswer = "<w:p $wNs $w14Ns " + `
    'w:rsidR="00000000" w:rsidDel="00000000" w:rsidP="00000000" ' + `
    'w:rsidRDefault="00000000" w:rsidRPr="00000000" w14:paraId="00000066">' + `
    "<w:pPr>$rPrBlock</w:pPr>" + `
    '<w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000">' + `
    "$runRPrBlock" + `
    '<w:t xml:space="preserve">' + $answerText + '</w:t>' + `
    '</w:r>' + `
    '</w:p>'
$pAnswer.Range.InsertXML($xmlAnswer)
